$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new "Status" value for row 8 (Investigate minified operators issue)
$ws.Range("F8").Value = "Resoved by SP"
$ws.Range("F8").WrapText = $true

# Add new row 14: a new issue under "Expansions" / "Medium" priority
$ws.Range("A14").Value = "Expansions"
$ws.Range("B14").Value = "Medium"
$ws.Range("C14").Value = "Return predicted reactions in order of rule that generated them"
$ws.Range("A14:C14").WrapText = $true
$ws.Rows.Item(14).RowHeight = 30

# Update the selected cell to match the saved view state
$ws.Range("C15").Select()
